# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F3 8519 -> 8520, F10 1214 -> 1215, F11 86 -> 87
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 8520
$wsExhibit.Range("F10").Value = 1215
$wsExhibit.Range("F11").Value = 87

# Sheet "全部类型" (sheet4): F3 8519 -> 8520, F14 1214 -> 1215, F15 86 -> 87
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 8520
$wsAll.Range("F14").Value = 1215
$wsAll.Range("F15").Value = 87
